$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two weekly date-groups (row2/row3 vs row4/row5) swap their
# Volumen/Precio values, effectively exchanging which week's data
# sits in rows 2-3 vs rows 4-5.

# Row 2 <- old Row 4 values
$ws.Range("D2").Value = 44223
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 3500
$ws.Range("O2").Value = 4000
$ws.Range("P2").Value = 3750
$ws.Range("S2").Value = 1875

# Row 3 <- old Row 5 values
$ws.Range("D3").Value = 44223
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 3000
$ws.Range("O3").Value = 3000
$ws.Range("P3").Value = 3000
$ws.Range("S3").Value = 1500

# Row 4 <- old Row 2 values
$ws.Range("D4").Value = 44559
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 6000
$ws.Range("O4").Value = 7000
$ws.Range("P4").Value = 6500
$ws.Range("S4").Value = 3250

# Row 5 <- old Row 3 values
$ws.Range("D5").Value = 44559
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 5000
$ws.Range("O5").Value = 5000
$ws.Range("P5").Value = 5000
$ws.Range("S5").Value = 2500

$wb.Save()
